$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold, border, centered) from the existing H1 header
# cell onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# PasteSpecial can clobber the values we just set (it pastes formats only,
# but reassign defensively in case the engine copies everything).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2-12
$values = @(
    @(7, 7),
    @(9, 9),
    @(6, 8),
    @(7, 8),
    @(8, 9),
    @(5, 7),
    @(9, 9),
    @(8, 9),
    @(8, 9),
    @(7, 8),
    @(7, 9)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
